$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 10:53"

# Polonia overtakes Japon in total cases -> swap labels so the sorted-by-
# total-cases layout stays intact; row 47 keeps the higher "Casos totales"
# number (Polonia, now updated) and row 48 keeps Japon's previous figures.
$ws.Range("A47").Value = "Polonia"
$ws.Range("A48").Value = "Japon"

# Row 25 (country rank 29): refreshed totals
$ws.Range("B25").Value = 228403
$ws.Range("C25").Value = 1987
$ws.Range("D25").Value = 159475
$ws.Range("E25").Value = 65240
$ws.Range("G25").Value = 65
$ws.Range("H25").Value = 3688

# Row 26 (country rank 30): refreshed totals
$ws.Range("B26").Value = 184268
$ws.Range("C26").Value = 3622
$ws.Range("D26").Value = 132055
$ws.Range("E26").Value = 44463
$ws.Range("G26").Value = 134
$ws.Range("H26").Value = 7750

# Row 28 (country rank 32): refreshed totals
$ws.Range("B28").Value = 128228
$ws.Range("C28").Value = 2430
$ws.Range("D28").Value = 59676
$ws.Range("E28").Value = 65842
$ws.Range("G28").Value = 54
$ws.Range("H28").Value = 2710

# Row 47: now Polonia's refreshed totals
$ws.Range("B47").Value = 69129
$ws.Range("C47").Value = 612
$ws.Range("D47").Value = 48593
$ws.Range("E47").Value = 18444
$ws.Range("G47").Value = 14
$ws.Range("H47").Value = 2092

# Row 48: now Japon's (previous) totals
$ws.Range("B48").Value = 69001
$ws.Range("D48").Value = 58428
$ws.Range("E48").Value = 9266
$ws.Range("H48").Value = 1307

# Row 52 (country rank 56): refreshed totals
$ws.Range("B52").Value = 56908
$ws.Range("C52").Value = 48
$ws.Range("E52").Value = 990

# Row 128 (country rank 132): refreshed totals
$ws.Range("B128").Value = 3102
$ws.Range("C128").Value = 1
$ws.Range("E128").Value = 207

# Row 131 (country rank 135): refreshed totals
$ws.Range("B131").Value = 2978
$ws.Range("C131").Value = 20
$ws.Range("D131").Value = 1901
$ws.Range("E131").Value = 991

# Row 137 (country rank 141): refreshed totals
$ws.Range("B137").Value = 2441
$ws.Range("C137").Value = 26
$ws.Range("D137").Value = 2150
$ws.Range("E137").Value = 227
